$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "55.833.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.41%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.451.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.38%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "484.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +10.13%  "

$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.505"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.48%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.456.17"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.68%  "

$ws.Range("E10").Value = "  +8.71%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0962"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.329"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.60%  "

$ws.Range("E13").Value = "  +1.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.876.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "55.867.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.34%  "

$ws.Range("E17").Value = "  +0.65%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.461.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.79%  "

$ws.Range("E19").Value = "  +6.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "315.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.29%  "

$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.36%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "58.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.22%  "

$ws.Range("E25").Value = "  +5.90%  "

$ws.Range("E26").Value = "  -1.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.159"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.570.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.89%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0775"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.20%  "

$ws.Range("E31").Value = "  -0.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "147.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.89%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.90%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.12"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.08%  "

$ws.Range("E36").Value = "  +8.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.69"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.85%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.848"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.41%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.996"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.27%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0547"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.71%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.597"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.12%  "

$ws.Range("E44").Value = "  +6.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "260.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +11.99%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0918"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.98%  "

$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.78%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +11.98%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0226"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.52%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.98%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.861.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.79%  "
